$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" - strikeouts) values for rows 2-6
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 3
